# Auto-generated edit script: update cryptos price/volume columns (D,E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.034.71"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "3.326.70"
$ws.Range("E3").Value = "  -4.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.96"
$ws.Range("E5").Value = "  -5.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "651.15"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.36"
$ws.Range("E7").Value = "  -12.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.413"
$ws.Range("E8").Value = "  -11.24%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.980"
$ws.Range("E10").Value = "  -11.56%  "
$ws.Range("D11").Value = "3.325.07"
$ws.Range("E11").Value = "  -4.22%  "
$ws.Range("E12").Value = "  -7.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.84"
$ws.Range("E13").Value = "  -7.20%  "
$ws.Range("D14").Value = "96.834.28"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("E15").Value = "  -4.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000250"
$ws.Range("E16").Value = "  -9.03%  "
$ws.Range("D17").Value = "3.941.87"
$ws.Range("E17").Value = "  -4.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.54"
$ws.Range("E18").Value = "  +3.08%  "
$ws.Range("D19").Value = "3.315.10"
$ws.Range("E19").Value = "  -4.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.70"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.479"
$ws.Range("E21").Value = "  +4.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.44"
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "490.54"
$ws.Range("E23").Value = "  -8.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.26"
$ws.Range("E24").Value = "  -10.37%  "
$ws.Range("E25").Value = "  -10.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.36"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "92.25"
$ws.Range("E27").Value = "  -10.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.95"
$ws.Range("E28").Value = "  -8.05%  "
$ws.Range("D29").Value = "3.496.12"
$ws.Range("E29").Value = "  -4.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -6.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.75"
$ws.Range("E32").Value = "  -6.55%  "
$ws.Range("E33").Value = "  -7.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.46"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.542"
$ws.Range("E36").Value = "  -8.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.84"
$ws.Range("E37").Value = "  -9.74%  "
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.48"
$ws.Range("E39").Value = "  -7.34%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -7.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "501.50"
$ws.Range("E42").Value = "  -7.62%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.67"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.821"
$ws.Range("E45").Value = "  -5.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0405"
$ws.Range("E46").Value = "  -8.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.32"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.40"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.35"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.11"
$ws.Range("E51").Value = "  -11.74%  "
